# Edit script for SaveSequenceDiagram.pptx
#
# Commit message: "Add ReadSequenceDiagram and introduced compareData method"
#
# The underlying diff removes a cluster of shapes from the (single) slide
# that made up the old "Logic"/SaveCommand sequence-diagram swimlane
# artwork, and nudges one remaining connector's geometry slightly so it
# still lines up with what is left behind.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Helper: find a shape on a slide by its (stable) Shape.Id, since shape
# names are not unique on this slide and indices shift as shapes are
# removed.
# ---------------------------------------------------------------------
function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Helper: Shape.Left/Top/Width/Height are single-precision points under
# PowerPoint's COM object model, so naively assigning `emu / 12700`
# can truncate to one EMU below the intended target once it round-trips
# back through the float32 storage. Binary-search the double-precision
# point value whose float32 cast maps (via floor) to exactly the EMU we
# want, so the saved OOXML matches byte-for-byte.
# ---------------------------------------------------------------------
function EmuToPt($emuTarget) {
    if ($emuTarget -eq 0) { return 0.0 }
    $lo = ($emuTarget / 12700.0)
    $hi = $lo + 0.001
    for ($i = 0; $i -lt 200; $i++) {
        $mid = ($lo + $hi) / 2.0
        $f = [float]$mid
        $back = [math]::Floor([double]$f * 12700.0)
        if ($back -ge $emuTarget) {
            $hi = $mid
        } else {
            $lo = $mid
        }
    }
    return $hi
}

# ---------------------------------------------------------------------
# Remove the shapes that belonged to the old "Logic" lane / SaveCommand
# sequence artwork (rectangles, connectors, labels and the grouped
# ":CommandResult" artwork) that the diff deletes wholesale.
# ---------------------------------------------------------------------
$idsToDelete = @(
    81,   # Rectangle 65 ("Logic" lane header)
    20,   # Straight Connector 19
    21,   # Rectangle 20
    28,   # Straight Arrow Connector 27
    29,   # TextBox 28 ("execute()")
    34,   # Straight Arrow Connector 33
    19,   # Rectangle 62 (":SaveCommand")
    66,   # Straight Arrow Connector 65
    69,   # Straight Arrow Connector 68
    65,   # Rectangle 64
    94,   # TextBox 93 ("X")
    223,  # Group 222 (":CommandResult" group)
    229,  # Straight Arrow Connector 228
    232,  # Straight Arrow Connector 231
    57    # TextBox 56 ("saveCommand()")
)

foreach ($id in $idsToDelete) {
    $shp = Get-ShapeById $s $id
    if ($shp -ne $null) {
        $shp.Delete()
    }
}

# ---------------------------------------------------------------------
# Nudge the remaining "Straight Arrow Connector 49" (id 50) so its
# geometry matches the post-edit layout now that neighbouring shapes are
# gone.
# ---------------------------------------------------------------------
$conn = Get-ShapeById $s 50
if ($conn -ne $null) {
    $conn.Left   = EmuToPt 4417379
    $conn.Top    = EmuToPt 3211220
    $conn.Width  = EmuToPt 2938302
    $conn.Height = EmuToPt 0
}
